$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '45.393.29'
$ws.Range('E2').Value = '  +6.46%  '

# Row 3
$ws.Range('D3').Value = '2.370.07'
$ws.Range('E3').Value = '  +2.98%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = '110.43'
$ws.Range('E5').Value = '  +5.02%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '312.99'
$ws.Range('E6').Value = '  +0.69%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +1.08%  '

# Row 8
$ws.Range('E8').Value = '  +0.10%  '

# Row 9
$ws.Range('D9').Value = '0.617'
$ws.Range('E9').Value = '  +2.22%  '

# Row 10
$ws.Range('D10').Value = '41.07'
$ws.Range('E10').Value = '  +3.35%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0920'
$ws.Range('E11').Value = '  +0.97%  '

# Row 12
$ws.Range('D12').Value = '8.55'
$ws.Range('E12').Value = '  +2.62%  '

# Row 13
$ws.Range('E13').Value = '  +1.67%  '

# Row 14
$ws.Range('D14').Value = '0.985'
$ws.Range('E14').Value = '  -0.21%  '

# Row 15
$ws.Range('D15').Value = '2.730.60'
$ws.Range('E15').Value = '  +3.00%  '

# Row 16
$ws.Range('D16').Value = '15.47'
$ws.Range('E16').Value = '  +1.47%  '

# Row 17
$ws.Range('D17').Value = '2.374.13'
$ws.Range('E17').Value = '  +3.20%  '

# Row 18
$ws.Range('D18').Value = '45.298.59'
$ws.Range('E18').Value = '  +5.84%  '

# Row 19
$ws.Range('D19').Value = '7.33'
$ws.Range('E19').Value = '  -0.48%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000107'
$ws.Range('E20').Value = '  +1.63%  '

# Row 21
$ws.Range('D21').Value = '13.12'
$ws.Range('E21').Value = '  -4.04%  '

# Row 22
$ws.Range('E22').Value = '  +0.20%  '

# Row 23
$ws.Range('D23').Value = '3.48'
$ws.Range('E23').Value = '  +0.56%  '

# Row 24
$ws.Range('D24').Value = '261.18'
$ws.Range('E24').Value = '  -2.62%  '

# Row 25
$ws.Range('E25').Value = '  +2.47%  '

# Row 26
$ws.Range('E26').Value = '  -0.41%  '

# Row 27
$ws.Range('D27').Value = '11.16'
$ws.Range('E27').Value = '  +2.18%  '

# Row 28
$ws.Range('D28').Value = '7.38'
$ws.Range('E28').Value = '  -6.74%  '

# Row 29
$ws.Range('D29').Value = '2.36'
$ws.Range('E29').Value = '  +2.37%  '

# Row 30
$ws.Range('D30').Value = '22.55'
$ws.Range('E30').Value = '  +1.64%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '38.20'
$ws.Range('E31').Value = '  +0.99%  '

# Row 32
$ws.Range('D32').Value = '0.0963'
$ws.Range('E32').Value = '  +11.59%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '169.80'
$ws.Range('E33').Value = '  +2.37%  '

# Row 34
$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  +5.01%  '

# Row 35
$ws.Range('D35').Value = '0.131'
$ws.Range('E35').Value = '  +0.05%  '

# Row 36
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '4.82'
$ws.Range('E36').Value = '  +4.00%  '

# Row 37
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  +3.29%  '

# Row 38
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.00'
$ws.Range('E38').Value = '  +8.21%  '

# Row 39
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '3.97'
$ws.Range('E39').Value = '  +10.15%  '

# Row 40
$ws.Range('D40').Value = '0.0357'
$ws.Range('E40').Value = '  -0.32%  '

# Row 41
$ws.Range('D41').Value = '1.73'
$ws.Range('E41').Value = '  +11.22%  '

# Row 42
$ws.Range('E42').Value = '  -4.72%  '

# Row 43
$ws.Range('D43').Value = '0.235'
$ws.Range('E43').Value = '  +2.60%  '

# Row 44
$ws.Range('D44').Value = '13.32'
$ws.Range('E44').Value = '  +8.81%  '

# Row 45
$ws.Range('D45').Value = '69.94'
$ws.Range('E45').Value = '  -1.52%  '

# Row 46
$ws.Range('E46').Value = '  -0.40%  '

# Row 47
$ws.Range('D47').Value = '81.99'
$ws.Range('E47').Value = '  +8.16%  '

# Row 48
$ws.Range('D48').Value = '113.27'
$ws.Range('E48').Value = '  +2.46%  '

# Row 49
$ws.Range('D49').Value = '9.35'
$ws.Range('E49').Value = '  +5.88%  '

# Row 50
$ws.Range('D50').Value = '5.55'
$ws.Range('E50').Value = '  +7.59%  '

# Row 51
$ws.Range('D51').Value = '1.640.30'
$ws.Range('E51').Value = '  -3.12%  '
